$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10/11 rework ----------------------------------------------------
# Old layout:
#   row10: B10            (blank,  style = "top border + fill-apply")
#   row11: B11=0 (style "fill-apply"), C11=14, D11=21, E11=4, F11=0
# New layout:
#   row10: B10=0 (style "fill-apply"), C10=14, D10=21, E10=4, F10=0
#   row11: B11             (blank,  style = "fill-apply", same family as B10)
$ws.Range("B10").Value = 0
$ws.Range("B10").Borders.LineStyle = -4142        # xlLineStyleNone - drop the old top border
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 0

$ws.Range("C11:F11").ClearContents()
$ws.Range("B11").ClearContents()                  # keeps its existing "fill-apply" style (s=10)

# --- Second (transposed) nonogram grid, rows 13-17 -----------------------
# 5x5 solution grid, same box-border look as the B5:F9 grid above
$ws.Range("B5:F9").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws.Application.CutCopyMode = $false

$grid = @(
  @(0,0,0,0,0),
  @(0,0,0,0,0),
  @(0,0,1,0,0),
  @(0,1,0,0,0),
  @(0,1,1,1,0)
)
$cols = @("B","C","D","E","F")
for ($r = 0; $r -lt 5; $r++) {
  $row = 13 + $r
  for ($c = 0; $c -lt 5; $c++) {
    $ws.Range($cols[$c] + $row).Value = $grid[$r][$c]
  }
}

# Row clue numbers in column A
$ws.Range("A13").Value = 1
$ws.Range("A14").Value = 2
$ws.Range("A15").Value = 3
$ws.Range("A16").Value = "2, 1"
$ws.Range("A17").Value = 4
$ws.Range("A16").HorizontalAlignment = -4152       # xlRight

# --- Selection / view state -----------------------------------------------
$ws.Range("D9").Select()
